$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 511, shifting existing rows 511-597 down to 512-598.
$ws.Rows(511).Insert()

# Populate the newly inserted row 511 with a new data point (same market/category
# metadata as the row that used to occupy position 511, but a new date and new
# price values).
$ws.Cells.Item(511, 1).Value2 = 4
$ws.Cells.Item(511, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(511, 3).Value2 = "Los Lagos"
$ws.Cells.Item(511, 4).Value2 = 45180
$ws.Cells.Item(511, 5).Value2 = 10
$ws.Cells.Item(511, 6).Value2 = 100112023
$ws.Cells.Item(511, 7).Value2 = "Brócoli"
$ws.Cells.Item(511, 8).Value2 = "Sin especificar"
$ws.Cells.Item(511, 9).Value2 = "Primera"
$ws.Cells.Item(511, 10).Value2 = 500
$ws.Cells.Item(511, 11).Value2 = 1400
$ws.Cells.Item(511, 12).Value2 = 1500
$ws.Cells.Item(511, 13).Value2 = 1450
$ws.Cells.Item(511, 14).Value2 = "$/unidad"
$ws.Cells.Item(511, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(511, 16).Value2 = 1450
$ws.Cells.Item(511, 17).Value2 = 1
$ws.Cells.Item(511, 18).Value2 = "Hortaliza"
